# Push Excel data 19th Aug
$wb = $excel.ActiveWorkbook

# --- addListItem sheet: bump the "fifteen" test user / code values ---
$wsAdd = $wb.Worksheets.Item("addListItem")
$wsAdd.Range("A2").Value = "UserfifteenM"
$wsAdd.Range("D2").Value = "ADLILC.8850"

# --- createUser sheet: bump the numeric increment used to build the test user ---
$wsCreate = $wb.Worksheets.Item("createUser")
$wsCreate.Range("A2").Value = 1053

# Leave a "clicked" selection on createUser (D10) before switching tabs,
# then land on addListItem (keeping its existing D2 selection) as the
# active sheet/tab for this save.
$wsCreate.Range("D10").Select()
$wsAdd.Activate()
$wsAdd.Range("D2").Select()
